$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "max" column (column C) entirely -- this shifts the
# "prediction" column (D -> C) and "rejection-f" column (E -> D) one
# position to the left, matching the new A1:D2 dimension.
$ws.Range("C1").EntireColumn.Delete()

# Update the remaining numeric prediction value in row 2.
$ws.Range("B2").Value = 23.43590354925201
